# Update "想去人数" (interest count) figures in column F across the four
# sheets of the workbook, reflecting a newer data pull (gh-pages output
# regenerated at commit 456a3b4).
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsAll = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibitions)
$wsExhibit.Range("F5").Value = 5640
$wsExhibit.Range("F8").Value = 2757
$wsExhibit.Range("F9").Value = 6409
$wsExhibit.Range("F10").Value = 177
$wsExhibit.Range("F11").Value = 1236
$wsExhibit.Range("F19").Value = 145
$wsExhibit.Range("F21").Value = 896
$wsExhibit.Range("F26").Value = 17
$wsExhibit.Range("F30").Value = 226
$wsExhibit.Range("F31").Value = 1151

# 演出 (Performances)
$wsShow.Range("F19").Value = 173
$wsShow.Range("F23").Value = 88
$wsShow.Range("F27").Value = 608

# 本地生活 (Local life)
$wsLocal.Range("F6").Value = 515

# 全部类型 (All types - aggregated view)
$wsAll.Range("F9").Value = 515
$wsAll.Range("F14").Value = 5640
$wsAll.Range("F17").Value = 2757
$wsAll.Range("F19").Value = 6409
$wsAll.Range("F21").Value = 177
$wsAll.Range("F22").Value = 1236
$wsAll.Range("F30").Value = 145
$wsAll.Range("F32").Value = 896
$wsAll.Range("F33").Value = 88
$wsAll.Range("F44").Value = 226
